$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 2.4
$ws.Range("R3").Value = 1.53
$ws.Range("W3").Value = 6
$ws.Range("AC3").Value = 7

# Row 4 updates
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.63
$ws.Range("Q4").Value = 2.4
$ws.Range("R4").Value = 1.53
$ws.Range("AM4").Value = 900

# Row 5 updates
$ws.Range("G5").Value = 2.05
$ws.Range("I5").Value = 3.5
$ws.Range("J5").Value = 2.88
$ws.Range("L5").Value = 4.33
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 1.73
$ws.Range("X5").Value = 9
$ws.Range("Z5").Value = 19
$ws.Range("AG5").Value = 8.5
$ws.Range("AH5").Value = 17
$ws.Range("AN5").Value = 4
$ws.Range("AO5").Value = 12
$ws.Range("AX5").Value = 21
$ws.Range("BA5").Value = 101

# Row 6 updates
$ws.Range("M6").Value = 1.07
$ws.Range("N6").Value = 9
$ws.Range("Q6").Value = 2.25
$ws.Range("R6").Value = 1.62
